$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2..49 down to 3..50)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new entry
$ws.Range("A2").Value = 2019
$ws.Range("B2").Value = "Andrés Gomez"

# Update selection to match the target state
$ws.Range("B3").Select()
